$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# --- Add the new "Old Town of Salamanca" coin row (row 32), mirroring the
# --- previous entry (row 31) for formatting, then filling in its own data ---
$ws.Range("A31:I31").Copy()
$ws.Range("A32:I32").PasteSpecial(-4122)

$ws.Range("A32").Value2 = 2025
$ws.Range("B32").Value2 = "Old Town of Salamanca"
$ws.Range("C32").Value2 = "Spanish UNESCO"
$ws.Range("D32").Value2 = "Obv: With mint logo - ""M"""
$ws.Range("E32").Value2 = "Rev: new map of Europe"
$ws.Range("F32").Value2 = ""
$ws.Range("G32").Value2 = "2.000.000"
$ws.Range("H32").Value2 = 1
$ws.Range("I32").Formula = '=IF(OR(AND(H32>1,H32<>"-")),"Can exchange","")'

# --- Re-create the H-column conditional formatting (duplicate-detector +
# --- color scale) for the new row, matching every other row in the column ---
$h32 = $ws.Range("H32")

$ct = $h32.FormatConditions.Add(9, 7, "*-")
$ct.Formula1 = '=NOT(ISERROR(SEARCH(("*-"),(H32))))'
$ct.Text = "*-"
$ct.Interior.Color = 16770459

$cs = $h32.FormatConditions.AddColorScale(3)
$cs.ColorScaleCriteria.Item(1).Type = 0
$cs.ColorScaleCriteria.Item(1).Value = 0
$cs.ColorScaleCriteria.Item(1).FormatColor.Color = 10526207
$cs.ColorScaleCriteria.Item(2).Type = 0
$cs.ColorScaleCriteria.Item(2).Value = 1
$cs.ColorScaleCriteria.Item(2).FormatColor.Color = 11716049
$cs.ColorScaleCriteria.Item(3).Type = 0
$cs.ColorScaleCriteria.Item(3).Value = 10
$cs.ColorScaleCriteria.Item(3).FormatColor.Color = 5287680

# --- Move the active selection to mirror where editing finished ---
$ws.Range("G36").Select()

$wb.Application.Calculate()
